$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before the current row 8 (pushing old rows 8-74 down
# to 9-75, matching dimension change A1:R74 -> A1:R75).
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with its data.
$ws.Range("A8").Value = 2
$ws.Range("B8").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = "2023-04-27"
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 100112032
$ws.Range("G8").Value = "Zapallo italiano"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 7000
$ws.Range("L8").Value = 7500
$ws.Range("M8").Value = 7250
$ws.Range("N8").Value = '$/caja 60 unidades'
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 121
$ws.Range("Q8").Value = 60
$ws.Range("R8").Value = "Hortaliza"
